$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C6 to a numeric value (20) instead of the "$" text placeholder
$ws.Range("C6").Value = 20

# Update the selection on the sheet to E1:F1048576 (entire columns E:F selection style)
$ws.Range("E1:F1048576").Select()
